$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 3311
$ws.Range('J3').Value = 3447
$ws.Range('E4').Value = 1992
$ws.Range('H4').Value = 1695
$ws.Range('J4').Value = 761
$ws.Range('J5').Value = 268
$ws.Range('J6').Value = 4066
$ws.Range('E7').Value = 25996
$ws.Range('H7').Value = 26005
$ws.Range('J7').Value = 11853

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 36
$ws.Range('J7').Value = 137

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 41
$ws.Range('J7').Value = 131

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 132
$ws.Range('J6').Value = 107
$ws.Range('J7').Value = 374

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J2').Value = 51
$ws.Range('J7').Value = 169

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 176
$ws.Range('J6').Value = 121
$ws.Range('J7').Value = 425

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J2').Value = 33
$ws.Range('J7').Value = 89

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 32
$ws.Range('J3').Value = 29
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 94
$ws.Range('J7').Value = 310

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J7').Value = 361
$ws.Range('J8').Value = 772
$ws.Range('J9').Value = 69
$ws.Range('J10').Value = 69
$ws.Range('J15').Value = 138
$ws.Range('J18').Value = 117
$ws.Range('J19').Value = 369
$ws.Range('J20').Value = 249
$ws.Range('J22').Value = 26
$ws.Range('J23').Value = 119
$ws.Range('J24').Value = 35
$ws.Range('J25').Value = 68
$ws.Range('J27').Value = 72
$ws.Range('J29').Value = 682
$ws.Range('J31').Value = 89
$ws.Range('J33').Value = 526
$ws.Range('J34').Value = 59
$ws.Range('J36').Value = 172
$ws.Range('J37').Value = 374
$ws.Range('J41').Value = 78
$ws.Range('J42').Value = 468
$ws.Range('J43').Value = 108
$ws.Range('J44').Value = 90
$ws.Range('J48').Value = 117
$ws.Range('J49').Value = 76
$ws.Range('J52').Value = 327
$ws.Range('J54').Value = 226
$ws.Range('J55').Value = 149
$ws.Range('J57').Value = 54
$ws.Range('J60').Value = 79
$ws.Range('E63').Value = 336
$ws.Range('H63').Value = 247
$ws.Range('J63').Value = 53
$ws.Range('J65').Value = 310
$ws.Range('J67').Value = 425
$ws.Range('J73').Value = 102
$ws.Range('J76').Value = 167
$ws.Range('J77').Value = 100
$ws.Range('J78').Value = 158
$ws.Range('J83').Value = 274
$ws.Range('J84').Value = 104
$ws.Range('J85').Value = 539
$ws.Range('J86').Value = 67
$ws.Range('J87').Value = 37
$ws.Range('J88').Value = 123
$ws.Range('J89').Value = 137
$ws.Range('J90').Value = 140
$ws.Range('J91').Value = 137
$ws.Range('J94').Value = 104
$ws.Range('J96').Value = 131
$ws.Range('J99').Value = 169
$ws.Range('J100').Value = 23
$ws.Range('E101').Value = 25996
$ws.Range('H101').Value = 26005
$ws.Range('J101').Value = 11853

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J6').Value = 77
$ws.Range('J7').Value = 274

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 143
$ws.Range('J3').Value = 168
$ws.Range('J5').Value = 20
$ws.Range('J6').Value = 169
$ws.Range('J7').Value = 526

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J3').Value = 14
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 60
$ws.Range('J7').Value = 226

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 209
$ws.Range('J3').Value = 234
$ws.Range('J6').Value = 172
$ws.Range('J7').Value = 682

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J5').Value = 17
$ws.Range('J6').Value = 141
$ws.Range('J7').Value = 369

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J6').Value = 26
$ws.Range('J7').Value = 90

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 21
$ws.Range('J7').Value = 117

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 34
$ws.Range('J7').Value = 167

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 133
$ws.Range('J7').Value = 539

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J6').Value = 39
$ws.Range('J7').Value = 78

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 105
$ws.Range('J7').Value = 468

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J2').Value = 20
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 40
$ws.Range('J3').Value = 52
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 158

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J3').Value = 31
$ws.Range('J7').Value = 149

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J4').Value = 8
$ws.Range('J7').Value = 35

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J6').Value = 26
$ws.Range('J7').Value = 119

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 63
$ws.Range('J6').Value = 22
$ws.Range('J7').Value = 137

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J3').Value = 79
$ws.Range('J7').Value = 249

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J6').Value = 62
$ws.Range('J7').Value = 117

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 63
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('J4').Value = 2
$ws.Range('J7').Value = 23

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 72
$ws.Range('J3').Value = 94
$ws.Range('J6').Value = 144
$ws.Range('J7').Value = 327

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 15
$ws.Range('J7').Value = 59

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 54
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 30
$ws.Range('J7').Value = 68

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J2').Value = 44
$ws.Range('J3').Value = 38
$ws.Range('J7').Value = 138

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J2').Value = 19
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J2').Value = 39
$ws.Range('J7').Value = 102

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J3').Value = 39
$ws.Range('J4').Value = 4
$ws.Range('J7').Value = 123

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 230
$ws.Range('J3').Value = 241
$ws.Range('J6').Value = 238
$ws.Range('J7').Value = 772

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 17
$ws.Range('J4').Value = 8
$ws.Range('J7').Value = 72

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J2').Value = 9
$ws.Range('J7').Value = 67

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 46
$ws.Range('J6').Value = 43
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J6').Value = 21
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J2').Value = 26
$ws.Range('J7').Value = 79

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J6').Value = 66
$ws.Range('J7').Value = 108

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J3').Value = 7
$ws.Range('J7').Value = 26

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J4').Value = 10
$ws.Range('J7').Value = 100

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 106
$ws.Range('J4').Value = 10
$ws.Range('J5').Value = 6
$ws.Range('J6').Value = 117
$ws.Range('J7').Value = 361

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('J6').Value = 20
$ws.Range('J7').Value = 37
